# Axar Patel (Delhi Capitals) batting log — "complete!! -> scrapping whole ipl"
# 1) Rename the sheet to the player's name
# 2) Insert a new leading "matchNo" column
# 3) Insert 5 new data rows (so the original single match row becomes one of six)
# 4) Re-populate every header + data cell for the finished A1:M7 table

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename sheet -------------------------------------------------------
$ws.Name = "Axar Patel"

# --- 2. Insert the new first column (shifts teamName..result from A:L -> B:M)
$ws.Columns("A").Insert()

# --- 3. Insert five rows above the existing data row 2 (old row 2 -> row 7)
$ws.Rows("2:6").Insert()

# --- 4. Force the whole table to be stored as text, matching the source data
#        (runs/balls/etc. are numeric-looking strings, not real numbers)
$ws.Range("A1:M7").NumberFormat = "@"

$headers = @("matchNo","teamName","batterName","states","runs","balls","fours","sixes","sr","opponentTeamName","venue","date","result")

for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws.Cells.Item(1, $c + 1).Value = $headers[$c]
}

$rows = @(
    @("Qualifier","Delhi Capitals","Axar Patel","","4","4","0","0","100.00","Kolkata Knight Riders","Sharjah","October 13","KKR won by 3 wickets (with 1 ball remaining)"),
    @("41st","Delhi Capitals","Axar Patel","c Ferguson b Iyer","0","5","0","0","0.00","Kolkata Knight Riders","Sharjah","September 28","KKR won by 3 wickets (with 10 balls remaining)"),
    @("50th","Delhi Capitals","Axar Patel","c Ali b Bravo","5","10","0","0","50.00","Chennai Super Kings","Dubai (DSC)","October 04","Capitals won by 3 wickets (with 2 balls remaining)"),
    @("Qualifier","Delhi Capitals","Axar Patel","c sub (MJ Santner) b Ali","10","11","1","0","90.90","Chennai Super Kings","Dubai (DSC)","October 10","Super Kings won by 4 wickets (with 2 balls remaining)"),
    @("36th","Delhi Capitals","Axar Patel","c Miller b Sakariya","12","7","0","1","171.42","Rajasthan Royals","Abu Dhabi","September 25","Capitals won by 33 runs"),
    @("46th","Delhi Capitals","Axar Patel","lbw b Boult","9","9","1","0","100.00","Mumbai Indians","Sharjah","October 02","Capitals won by 4 wickets (with 5 balls remaining)")
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $row = $rows[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}
